# Auto-generated edit script: updates cryptos list (crypto prices/volumes) per commit
# "Updated cryptos list on Sat Sep 28 02:52:34 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.140.66"
$ws.Range("E2").Value = "  +1.81%  "

# Row 3
$ws.Range("D3").Value = "2.700.51"
$ws.Range("E3").Value = "  +3.00%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.32"
$ws.Range("E5").Value = "  +1.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.95"
$ws.Range("E6").Value = "  +2.94%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +1.12%  "

# Row 9
$ws.Range("D9").Value = "0.126"
$ws.Range("E9").Value = "  +7.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.03"
$ws.Range("E10").Value = "  +4.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.405"
$ws.Range("E11").Value = "  +0.00%  "

# Row 12
$ws.Range("E12").Value = "  +1.05%  "

# Row 13
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  +13.05%  "

# Row 14
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "30.58"
$ws.Range("E14").Value = "  +5.18%  "

# Row 15
$ws.Range("D15").Value = "3.190.62"
$ws.Range("E15").Value = "  +3.22%  "

# Row 16
$ws.Range("D16").Value = "65.991.72"
$ws.Range("E16").Value = "  +1.69%  "

# Row 17
$ws.Range("D17").Value = "2.701.30"
$ws.Range("E17").Value = "  +2.97%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.79"
$ws.Range("E18").Value = "  +2.53%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.92"
$ws.Range("E19").Value = "  +1.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.84"
$ws.Range("E20").Value = "  +7.14%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "360.26"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.41"
$ws.Range("E22").Value = "  +4.29%  "

# Row 23
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000114"
$ws.Range("E24").Value = "  +20.78%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.94"
$ws.Range("E25").Value = "  +6.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.66"
$ws.Range("E26").Value = "  +0.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.69"
$ws.Range("E27").Value = "  +2.42%  "

# Row 28
$ws.Range("E28").Value = "  +4.80%  "

# Row 29
$ws.Range("E29").Value = "  +3.56%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +4.87%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "543.97"
$ws.Range("E31").Value = "  +3.20%  "

# Row 32
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("E33").Value = "  +1.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.74"
$ws.Range("E34").Value = "  +6.88%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.51"
$ws.Range("E35").Value = "  +0.66%  "

# Row 36
$ws.Range("D36").Value = "0.437"
$ws.Range("E36").Value = "  +2.88%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.90"
$ws.Range("E37").Value = "  +2.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "163.90"
$ws.Range("E38").Value = "  +0.84%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.01"
$ws.Range("E39").Value = "  +0.16%  "

# Row 40
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.01%  "

# Row 41
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "170.97"
$ws.Range("E41").Value = "  +4.02%  "

# Row 42
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.01%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.52"
$ws.Range("E43").Value = "  +1.17%  "

# Row 44
$ws.Range("E44").Value = "  +2.04%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.37"
$ws.Range("E45").Value = "  +7.85%  "

# Row 46
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0627"
$ws.Range("E46").Value = "  +2.56%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "23.84"
$ws.Range("E47").Value = "  +3.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0269"
$ws.Range("E48").Value = "  +3.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.662"
$ws.Range("E49").Value = "  +1.71%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.90"
$ws.Range("E50").Value = "  +7.68%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0995"
$ws.Range("E51").Value = "  +1.81%  "
